$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5650793650793651
$ws.Range("J2").Value = 0.01904761904761905
$ws.Range("P2").Value = 0.1333333333333333
$ws.Range("S2").Value = 0.09206349206349207
$ws.Range("B3").Value = 0.01075268817204301
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.02150537634408602
$ws.Range("P3").Value = 0.7258064516129032
$ws.Range("S3").Value = 0.2096774193548387
$ws.Range("J4").Value = 0.1481481481481481
$ws.Range("P4").Value = 0.7592592592592593
$ws.Range("S4").Value = 0.09259259259259259
$ws.Range("B6").Value = 0.06274509803921569
$ws.Range("D6").Value = 0.01568627450980392
$ws.Range("F6").Value = 0.07450980392156863
$ws.Range("J6").Value = 0.2392156862745098
$ws.Range("O6").Value = 0.02745098039215686
$ws.Range("Q6").Value = 0.1607843137254902
$ws.Range("R6").Value = 0.09411764705882353
$ws.Range("S6").Value = 0.3254901960784314
$ws.Range("B7").Value = 0.05851063829787234
$ws.Range("D7").Value = 0.0425531914893617
$ws.Range("F7").Value = 0.03191489361702127
$ws.Range("J7").Value = 0.2180851063829787
$ws.Range("O7").Value = 0.01595744680851064
$ws.Range("Q7").Value = 0.1170212765957447
$ws.Range("R7").Value = 0.06382978723404255
$ws.Range("S7").Value = 0.4521276595744681
$ws.Range("B8").Value = 0.1014492753623188
$ws.Range("D8").Value = 0.01863354037267081
$ws.Range("F8").Value = 0.06004140786749482
$ws.Range("J8").Value = 0.1118012422360248
$ws.Range("O8").Value = 0.02070393374741201
$ws.Range("Q8").Value = 0.1842650103519669
$ws.Range("R8").Value = 0.1076604554865424
$ws.Range("S8").Value = 0.3954451345755693
$ws.Range("B9").Value = 0.1022222222222222
$ws.Range("D9").Value = 0.02222222222222222
$ws.Range("F9").Value = 0.04888888888888889
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.02666666666666667
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.1066666666666667
$ws.Range("S9").Value = 0.3822222222222222
$ws.Range("B10").Value = 0.1081646894626657
$ws.Range("D10").Value = 0.0209351011863224
$ws.Range("F10").Value = 0.07676203768318214
$ws.Range("J10").Value = 0.1311933007676204
$ws.Range("O10").Value = 0.02372644801116539
$ws.Range("Q10").Value = 0.2107466852756455
$ws.Range("R10").Value = 0.08094905792044661
$ws.Range("S10").Value = 0.3475226796929519
$ws.Range("G11").Value = 0.163265306122449
$ws.Range("J11").Value = 0.1282798833819242
$ws.Range("K11").Value = 0.2419825072886297
$ws.Range("L11").Value = 0.434402332361516
$ws.Range("S11").Value = 0.03206997084548105
$ws.Range("G12").Value = 0.6923076923076923
$ws.Range("J12").Value = 0.2115384615384615
$ws.Range("K12").Value = 0.01923076923076923
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.05128205128205128
$ws.Range("G13").Value = 0.5918367346938775
$ws.Range("J13").Value = 0.3469387755102041
$ws.Range("S13").Value = 0.06122448979591837
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.01520912547528517
$ws.Range("H15").Value = 0.155893536121673
$ws.Range("I15").Value = 0.05703422053231939
$ws.Range("J15").Value = 0.2927756653992395
$ws.Range("K15").Value = 0.07604562737642585
$ws.Range("M15").Value = 0.007604562737642586
$ws.Range("O15").Value = 0.08365019011406843
$ws.Range("S15").Value = 0.311787072243346
$ws.Range("F16").Value = 0.004739336492890996
$ws.Range("H16").Value = 0.1421800947867299
$ws.Range("I16").Value = 0.08056872037914692
$ws.Range("J16").Value = 0.4549763033175355
$ws.Range("K16").Value = 0.1184834123222749
$ws.Range("M16").Value = 0.004739336492890996
$ws.Range("O16").Value = 0.07582938388625593
$ws.Range("S16").Value = 0.1184834123222749
$ws.Range("F17").Value = 0.02191235059760956
$ws.Range("H17").Value = 0.1772908366533865
$ws.Range("I17").Value = 0.08167330677290836
$ws.Range("J17").Value = 0.3864541832669323
$ws.Range("K17").Value = 0.08764940239043825
$ws.Range("M17").Value = 0.0199203187250996
$ws.Range("N17").Value = 0.00398406374501992
$ws.Range("O17").Value = 0.0697211155378486
$ws.Range("S17").Value = 0.151394422310757
$ws.Range("F18").Value = 0.03083700440528634
$ws.Range("H18").Value = 0.1850220264317181
$ws.Range("I18").Value = 0.105726872246696
$ws.Range("J18").Value = 0.4140969162995595
$ws.Range("K18").Value = 0.08370044052863436
$ws.Range("M18").Value = 0.03083700440528634
$ws.Range("O18").Value = 0.05726872246696035
$ws.Range("S18").Value = 0.09251101321585903
$ws.Range("F19").Value = 0.01626591230551627
$ws.Range("H19").Value = 0.1944837340876945
$ws.Range("I19").Value = 0.09193776520509193
$ws.Range("J19").Value = 0.3691654879773691
$ws.Range("K19").Value = 0.103960396039604
$ws.Range("M19").Value = 0.0198019801980198
$ws.Range("N19").Value = 0.002121640735502122
$ws.Range("O19").Value = 0.06647807637906648
$ws.Range("S19").Value = 0.1357850070721358
